$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.300.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -5.32%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.351.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.06%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'556.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'170.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.44%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.149"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'55.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.0000263"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'8.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.20%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'3.918.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.31%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.339.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.47%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'17.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.50%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'63.369.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.28%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'11.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.969"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'399.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.71%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'4.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.28%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'81.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.08%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'13.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.28%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'10.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.79%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.55%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'8.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.38%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'29.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.31%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'6.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'584.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.18%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'11.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.105"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.38%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'58.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.148"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'35.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.95%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "Stacks"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'3.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "Maker"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'3.163.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.365"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "  -11.08%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'2.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'2.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.61%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'3.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.71%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.0399"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.49%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "  -6.04%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.127"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'134.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.13%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'8.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.15%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'2.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.86%  "
$ws.Range("E51").Style = "Normal"
